$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "expectedParkingTime = 0"
$ws.Range("E2").Value = "expectedParkingTime = 0"
$ws.Range("C4").Value = "expectedParkingTime = 0"
$ws.Range("E4").Value = "expectedParkingTime = 0"
